$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 0.0055
    3  = 0.0065
    4  = 0.005
    5  = 0.0055
    6  = 0.0065
    7  = 0.0065
    8  = 0.006
    9  = 0.006
    10 = 0.0065
    11 = 0.0055
    12 = 0.006
    13 = 0.0065
    14 = 0.0055
    15 = 0.0055
    16 = 0.0045
    17 = 0.0055
    18 = 0.005
    19 = 0.006
    20 = 0.0045
    21 = 0.004
    22 = 0.0045
    23 = 0.004
    24 = 0.0035
    25 = 0.002
    26 = 0.002
    27 = 0.0015
    28 = 0.001
    29 = 0.004
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 10).Value = $values[$row]
}

$ws.Range("O26").Select()
